# Daily attendance processing - 2025-12-11 18:36:15
# Reorder the "Recorded By" (column G) value on the Session Analysis
# Results sheet: when the list of recorders starts with "System, ",
# move that leading "System" token to the end of the comma-separated list.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#      "System, backup@backdoor.com, system" -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -ne $val -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value = "$rest, System"
    }
}
